$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before N (currently "hora"), shifting "hora" and "Id Ugm" right.
$ws.Range("N1").EntireColumn.Insert()

# New header cell for the inserted column.
$ws.Range("N1").Value = "universidad"

# Remove the now-stale placeholder empty rows (2-9); only the header row remains used.
$ws.Range("A2:A9").EntireRow.Delete()

# Match the post-edit selection state (column N selected, active cell N1).
$ws.Range("N1").EntireColumn.Select()
